# Consumer.xlsx update: add new tickers across sheets + rework the
# "Leisure" sheet (formerly the odd-one-out "Main" layout) to match the
# standard Name/Ticker/Price/MC layout used by the other sheets.

$wb = $excel.ActiveWorkbook

$shFoodBev = $wb.Worksheets.Item("Food-Beverages")
$shApparel = $wb.Worksheets.Item("Apparel")
$shRetail  = $wb.Worksheets.Item("Retail")
$shRest    = $wb.Worksheets.Item("Restaurants")
$shLeisure = $wb.Worksheets.Item("Leisure")
$shNonDur  = $wb.Worksheets.Item("NonDurable")

# --- New data rows, entered in the same order the shared-string table
# --- picked them up (interleaved across sheets): this keeps the
# --- sharedStrings.xml unique-string ordering identical to the target.

# Retail!23  Woolworths / WOW AU
$shRetail.Range("A23").Value = "x"
$shRetail.Range("B23").Value = "Woolworths"
$shRetail.Range("C23").Value = "WOW AU"

# Apparel!13  Adidas / ADS GR
$shApparel.Range("A13").Value = "x"
$shApparel.Range("B13").Value = "Adidas"
$shApparel.Range("C13").Value = "ADS GR"

# Food-Beverages!30  Givaudan / GIVN SW
$shFoodBev.Range("A30").Value = "x"
$shFoodBev.Range("B30").Value = "Givaudan"
$shFoodBev.Range("C30").Value = "GIVN SW"

# Retail!24  Loblaws / L CN
$shRetail.Range("A24").Value = "x"
$shRetail.Range("B24").Value = "Loblaws"
$shRetail.Range("C24").Value = "L CN"

# Food-Beverages!31  Tyson Foods / TSN
$shFoodBev.Range("A31").Value = "x"
$shFoodBev.Range("B31").Value = "Tyson Foods"
$shFoodBev.Range("C31").Value = "TSN"

# Food-Beverages!32  International Flavor / IFF
$shFoodBev.Range("A32").Value = "x"
$shFoodBev.Range("B32").Value = "International Flavor"
$shFoodBev.Range("C32").Value = "IFF"

# Retail!25  Ahold / AD NA
$shRetail.Range("A25").Value = "x"
$shRetail.Range("B25").Value = "Ahold"
$shRetail.Range("C25").Value = "AD NA"

# Retail!26  Ross Stores / ROST
$shRetail.Range("A26").Value = "x"
$shRetail.Range("B26").Value = "Ross Stores"
$shRetail.Range("C26").Value = "ROST"

# Leisure!5  Las Vegas Sands / LVS  (sheet is being reshaped below, before
# we actually type these two rows in - the header/reshape happens first so
# the row numbers line up, then the cells get their values)
# -- reshape happens just below, then rows 5-6 get filled in.

# --- Rework the Leisure sheet layout -----------------------------------
# Before: a bespoke "Main" sheet (A1 title, B/C only, no header row, no
# frozen panes). After: same Name/Ticker/Price/MC header row + frozen
# panes as every other sheet, plus two more holdings appended.

# the two existing holdings (previously B2:C2 / B3:C3, with no header
# row) shift down one row to make way for the new header row
$oldRow2Name = $shLeisure.Range("B2").Value()
$oldRow2Ticker = $shLeisure.Range("C2").Value()
$oldRow3Name = $shLeisure.Range("B3").Value()
$oldRow3Ticker = $shLeisure.Range("C3").Value()

$shLeisure.Range("B2").Value = "Name"
$shLeisure.Range("C2").Value = "Ticker"
$shLeisure.Range("D2").Value = "Price"
$shLeisure.Range("E2").Value = "MC"

$shLeisure.Range("A3").Value = "x"
$shLeisure.Range("B3").Value = $oldRow2Name
$shLeisure.Range("C3").Value = $oldRow2Ticker

$shLeisure.Range("A4").Value = "x"
$shLeisure.Range("B4").Value = $oldRow3Name
$shLeisure.Range("C4").Value = $oldRow3Ticker

# Las Vegas Sands / LVS
$shLeisure.Range("A5").Value = "x"
$shLeisure.Range("B5").Value = "Las Vegas Sands"
$shLeisure.Range("C5").Value = "LVS"

# Food-Beverages!33  Hormel Foods / HRL
$shFoodBev.Range("A33").Value = "x"
$shFoodBev.Range("B33").Value = "Hormel Foods"
$shFoodBev.Range("C33").Value = "HRL"

# Leisure!6  Galaxy Entertainment / 27 HK
$shLeisure.Range("A6").Value = "x"
$shLeisure.Range("B6").Value = "Galaxy Entertainment"
$shLeisure.Range("C6").Value = "27 HK"

# NonDurable!14  Beiersdorf / BEI GR
$shNonDur.Range("A14").Value = "x"
$shNonDur.Range("B14").Value = "Beiersdorf"
$shNonDur.Range("C14").Value = "BEI GR"

# --- Leisure sheet cosmetics: column widths + frozen header like the
# --- rest of the workbook.
$shLeisure.Range("C3").Select()
$shLeisure.Columns.Item(1).ColumnWidth = 4.166666666666667
$shLeisure.Columns.Item(2).ColumnWidth = 18.3
$excel.ActiveWindow.FreezePanes = $true

# --- Final selections / active sheet, in the order the user last
# --- touched each sheet (NonDurable ends up active/tabSelected, matching
# --- the saved workbook).
$shFoodBev.Range("D33").Select()
$shApparel.Range("B14").Select()
$shRetail.Range("B27").Select()
$shRest.Range("B2:E2").Select()
$shLeisure.Range("B7").Select()
$shNonDur.Range("B15").Select()

$shNonDur.Activate()
